# Daily attendance processing - 2026-01-13 10:38:35
#
# Normalizes the "Recorded By" (column G) values: any occurrence of the
# exact token "System" within the comma-separated list of recorders is
# moved to the end of the list, while the relative order of the other
# tokens is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ", "
    $others = @()
    $systemCount = 0

    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemCount = $systemCount + 1
        } else {
            $others += $p
        }
    }

    if ($systemCount -eq 0) { continue }

    $newParts = $others
    for ($i = 0; $i -lt $systemCount; $i++) {
        $newParts += "System"
    }
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
